# Update the answer values in the "two-digit number divided by one-digit
# number" practice table. The worksheet has one table; every 4th row
# (rows 1, 5, 9, 13, 17) holds five division problems, the rows in
# between are left blank for the student's work. We overwrite the text
# of each of those 25 cells with the new values from the latest
# generated answer set.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of row -> list of new cell values (columns 1-5)
$updates = @{
    1  = @("68÷6=11, 2", "68÷2=34, 0", "81÷8=10, 1", "76÷5=15, 1", "67÷6=11, 1")
    5  = @("21÷2=10, 1", "61÷5=12, 1", "88÷4=22, 0", "62÷8=7, 6",  "80÷5=16, 0")
    9  = @("74÷7=10, 4", "88÷4=22, 0", "64÷3=21, 1", "97÷8=12, 1", "84÷9=9, 3")
    13 = @("54÷3=18, 0", "16÷8=2, 0",  "68÷5=13, 3", "92÷3=30, 2", "79÷3=26, 1")
    17 = @("12÷8=1, 4",  "28÷4=7, 0",  "79÷3=26, 1", "90÷4=22, 2", "39÷5=7, 4")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
